$wb = $excel.ActiveWorkbook

# Template sheet to copy from
$greece = $wb.Worksheets.Item("Greece")

# --- Netherlands ---
$greece.Copy($null, $greece)
$netherlands = $wb.Worksheets.Item($wb.Worksheets.Count)
$netherlands.Name = "Netherlands"
$netherlands.Range("B4").Select()
$netherlands.Range("B4").Value = "NGC-3144/T2202"
$netherlands.Range("B2").Value = "Netherlands Market"

# --- Austria ---
$netherlands.Copy($null, $netherlands)
$austria = $wb.Worksheets.Item($wb.Worksheets.Count)
$austria.Name = "Austria"
$austria.Range("B4").Select()
$austria.Range("B4").Value = "NGC-3817/T2309"
$austria.Range("B2").Value = "Austria Market"

# --- Denmark ---
$austria.Copy($null, $austria)
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B4").Select()
$denmark.Range("B4").Value = "NGC-2913/T2790"
$denmark.Range("B2").Value = "Denmark Market"

# The Netherlands tab ends up being the active one (per target workbook state)
$netherlands.Activate()
$netherlands.Range("B4").Select()
